# Adds columns I ("I0") and J ("IF") to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1) - copy the formatting from the existing "IP" header cell (H1)
# so the new headers share the same (bold, bordered, centered) style, then set text.
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(1, 10).Value = "IF"

# Data rows 2..69: (row, I value, J value)
$data = @(
    @(2, 4, 5),
    @(3, 6, 6),
    @(4, 9, 9),
    @(5, 7, 7),
    @(6, 7, 7),
    @(7, 7, 7),
    @(8, 9, 9),
    @(9, 10, 10),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 9, 9),
    @(13, 8, 8),
    @(14, 1, 1),
    @(15, 1, 2),
    @(16, 1, 2),
    @(17, 9, 9),
    @(18, 1, 1),
    @(19, 10, 10),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 7, 7),
    @(26, 7, 7),
    @(27, 9, 9),
    @(28, 7, 7),
    @(29, 6, 6),
    @(30, 7, 7),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 7, 7),
    @(34, 6, 6),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 7, 7),
    @(38, 9, 9),
    @(39, 7, 7),
    @(40, 7, 7),
    @(41, 6, 6),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 7, 7),
    @(45, 9, 9),
    @(46, 7, 7),
    @(47, 7, 7),
    @(48, 8, 8),
    @(49, 6, 6),
    @(50, 12, 12),
    @(51, 8, 8),
    @(52, 4, 4),
    @(53, 7, 7),
    @(54, 7, 7),
    @(55, 4, 4),
    @(56, 6, 6),
    @(57, 9, 9),
    @(58, 7, 7),
    @(59, 3, 4),
    @(60, 9, 9),
    @(61, 9, 9),
    @(62, 7, 7),
    @(63, 6, 6),
    @(64, 7, 7),
    @(65, 6, 6),
    @(66, 7, 7),
    @(67, 6, 6),
    @(68, 6, 6),
    @(69, 7, 7)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
